$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 684-686, shifting existing rows 684-701 down to 687-704
$ws.Rows("684:686").Insert()

# Row 684
$ws.Cells.Item(684, 1).Value = 10
$ws.Cells.Item(684, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(684, 3).Value = "La Araucanía"
$ws.Cells.Item(684, 4).Value = 45239
$ws.Cells.Item(684, 5).Value = 9
$ws.Cells.Item(684, 6).Value = 100114014
$ws.Cells.Item(684, 7).Value = "Betarraga"
$ws.Cells.Item(684, 8).Value = "Sin especificar"
$ws.Cells.Item(684, 9).Value = "Primera"
$ws.Cells.Item(684, 10).Value = 50
$ws.Cells.Item(684, 11).Value = 10000
$ws.Cells.Item(684, 12).Value = 10000
$ws.Cells.Item(684, 13).Value = 10000
$ws.Cells.Item(684, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(684, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(684, 16).Value = 833
$ws.Cells.Item(684, 17).Value = 12
$ws.Cells.Item(684, 18).Value = "Hortaliza"

# Row 685
$ws.Cells.Item(685, 1).Value = 10
$ws.Cells.Item(685, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(685, 3).Value = "La Araucanía"
$ws.Cells.Item(685, 4).Value = 45239
$ws.Cells.Item(685, 5).Value = 9
$ws.Cells.Item(685, 6).Value = 100114014
$ws.Cells.Item(685, 7).Value = "Betarraga"
$ws.Cells.Item(685, 8).Value = "Sin especificar"
$ws.Cells.Item(685, 9).Value = "Primera"
$ws.Cells.Item(685, 10).Value = 2000
$ws.Cells.Item(685, 11).Value = 800
$ws.Cells.Item(685, 12).Value = 800
$ws.Cells.Item(685, 13).Value = 800
$ws.Cells.Item(685, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(685, 15).Value = "Región Metropolitana"
$ws.Cells.Item(685, 16).Value = 160
$ws.Cells.Item(685, 17).Value = 5
$ws.Cells.Item(685, 18).Value = "Hortaliza"

# Row 686
$ws.Cells.Item(686, 1).Value = 10
$ws.Cells.Item(686, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(686, 3).Value = "La Araucanía"
$ws.Cells.Item(686, 4).Value = 45239
$ws.Cells.Item(686, 5).Value = 9
$ws.Cells.Item(686, 6).Value = 100114014
$ws.Cells.Item(686, 7).Value = "Betarraga"
$ws.Cells.Item(686, 8).Value = "Sin especificar"
$ws.Cells.Item(686, 9).Value = "Primera"
$ws.Cells.Item(686, 10).Value = 1000
$ws.Cells.Item(686, 11).Value = 800
$ws.Cells.Item(686, 12).Value = 800
$ws.Cells.Item(686, 13).Value = 800
$ws.Cells.Item(686, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(686, 15).Value = "Región del Maule"
$ws.Cells.Item(686, 16).Value = 160
$ws.Cells.Item(686, 17).Value = 5
$ws.Cells.Item(686, 18).Value = "Hortaliza"
